$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Global font rename: TimesNewToman -> Times New Roman (applies everywhere)
# ---------------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Font.Name = "TimesNewToman"
$find.Replacement.ClearFormatting()
$find.Replacement.Font.Name = "Times New Roman"
$find.Execute("", $false, $false, $false, $false, $false, $true, 1, $true, "", 2)

# ---------------------------------------------------------------------------
# 2) Simple text replacements (title, author, username, email, sentences)
# ---------------------------------------------------------------------------
function Replace-Text($oldText, $newText) {
    $rng = $d.Content
    $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
}

Replace-Text "Celestial Symphony: The Rhythms of the Universe" "The Allure of History: A Journey Through Time"
Replace-Text "Isabella Maxwell" "Alex Watson"
Replace-Text "isabella" "alex"
Replace-Text "maxwell@astronomycenter" "watson@eduworld"

Replace-Text "Within the vast canvas of the cosmic tapestry, celestial bodies engage in an intricate dance, governed by the harmonious laws of gravitation" "History beckons us, like an alluring whisper from times gone by"

Replace-Text " From the gentle waltz of our solar system's planets to the whirling dervishes of distant galaxies, the cosmos pulsates with an unseen symphony" " It is a kaleidoscope of human experiences, triumphs and follies, wisdom and folly, painted on the canvas of centuries"

Replace-Text " This symphony is a symphony of motion, dictated by the gravitational forces that bind celestial bodies together" " With each chapter, it holds a mirror to our present, shedding light on our origins, shaping our identities, and guiding us into the future"

Replace-Text "It is a symphony of time, measured by the ebb and flow of stellar ages, the rise and fall of civilizations, and the birth and death of stars" "In its vast expanse, history unveils a tapestry woven with countless threads, each representing the lives of individuals who have shaped our world"

Replace-Text " Every celestial object, from the smallest comet to the grandest supermassive black hole, plays a role in this cosmic orchestra, contributing to the intricate composition that orchestrates the universe" " Like intricate brushstrokes, their actions, decisions, and sacrifices add color and texture, creating a vibrant masterpiece that reveals the human capacity for both great achievements and devastating failures"

Replace-Text "The study of these celestial motions, known as celestial mechanics, delves into the underlying principles that govern the dynamics of the universe" "Furthermore, history teaches us the art of empathy and perspective"

Replace-Text " Scientists, like maestros of the universe, analyze the ballet of planets, the pirouette of stars, and the majestic procession of galaxies" " As we journey through the annals of time, we encounter diverse cultures, beliefs, and ways of life"

Replace-Text " Through this meticulous examination, they unravel the mysteries of the cosmos, revealing its hidden harmonies and unlocking its secrets" " We learn to appreciate the richness of human existence, and we begin to understand why people think, feel, and act as they do"

Replace-Text "The symphony of the universe unfolds through the graceful dance of celestial bodies, governed by the gravitational forces that orchestrate the cosmos" "History is a tapestry of human experiences, unveiling the richness and complexity of our shared past"

Replace-Text " Celestial mechanics, like a conductor of the universal orchestra, analyzes this dance, deciphering the principles that govern its rhythm" " It illuminates our present, shaping our identities, and guiding us into the future"

Replace-Text " The study of this cosmic choreography grants us insights into the workings of the universe, allowing us to appreciate the profound beauty and intricate interconnectedness of all things celestial" " Through its stories of triumphs and follies, wisdom and folly, history teaches us empathy, perspective, and the profound interconnectedness of humanity"

# ---------------------------------------------------------------------------
# 3) Insert brand-new sentences that did not exist before
# ---------------------------------------------------------------------------

# (a) After "...devastating failures." insert a new sentence, before the <w:br/>
$rng = $d.Content
$rng.Find.Execute("devastating failures.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.InsertAfter(" Amidst the ebb and flow of civilizations and empires, we discover the timeless struggles of humanity - the quest for power, the pursuit of justice, the yearning for freedom.")

# (b) After "...act as they do" insert two new sentences, before the trailing "."
$rng = $d.Content
$rng.Find.Execute("act as they do", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.InsertAfter(". This understanding fosters tolerance, compassion, and the realization that we are all part of a shared human story")

# (c) In the Summary paragraph, after "...interconnectedness of humanity" insert two new sentences
$rng = $d.Content
$rng.Find.Execute("interconnectedness of humanity", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.InsertAfter(". It is a beacon that sheds light on our origins, offering invaluable insights into who we are and how we can navigate the challenges of an ever-changing world")

# ---------------------------------------------------------------------------
# 4) Append a new empty paragraph at the very end of the document
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$endRng = $lastPara.Range
$endRng.Collapse(0)
$endRng.InsertBefore("`r")
